$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("test_schedule_1")

# Update file path strings to relative paths (remove the absolute prefix)
$ws1.Range("C2").Value = "test_files\vids\test.mp4"
$ws1.Range("C3").Value = "test_files\vids\test2.mp4"
$ws1.Range("C4").Value = "test_files\vids\test4.mp4"

# Update the active selection on the sheet from F7 to D11
$ws1.Range("D11").Select()
